$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Update existing rows 2-5: Befundtext (G) becomes "Vorname Name", E-Nummer (H) gets new fall numbers ---
$ws.Range("G2").Value = "Stephan Frank"
$ws.Range("H2").Value = "A/2001/200592"

$ws.Range("G3").Value = "Klaudis Apfel"
$ws.Range("H3").Value = "A/2002/200591"

$ws.Range("G4").Value = "Erika Mustermann"
$ws.Range("H4").Value = "A/2003/200391"

$ws.Range("G5").Value = "Max Mustermann"
$ws.Range("H5").Value = "A/2004/200591"

# --- Update row 6: Alternativer Name (C) changes, Befundtext (G) becomes "Klaus Kleber Excel 1" ---
$ws.Range("C6").Value = "TEST 1 update"
$ws.Range("G6").Value = "Klaus Kleber Excel 1"

# --- Add new row 7 (another "Fall" for Klaus Kleber) by copying row 6's formatting down ---
$ws.Range("A6:I6").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A7").Value = "Klaus"
$ws.Range("B7").Value = "Kleber"
$ws.Range("C7").Value = "TEST 2 kein update"
$ws.Range("D7").Value = "Klassenweg"
$ws.Range("E7").Value2 = $ws.Range("E6").Value2
$ws.Range("F7").Value = 99999
$ws.Range("G7").Value = "Klaus Kleber Excel 2"
$ws.Range("H7").Value = "A/2000/123456"
$ws.Range("I7").Value = "Nebenbefund"

$ws.Range("G7").Select()
